# "10Th - MB for single stock and added new group"
#
# The MarketBeat single-stock analyst-rating grid gets a new reporting
# group: two duplicate "Jun_26" columns (for a rating change that needs to
# stand out) plus a brand-new "Jun_27" column, all inserted to the left of
# the existing Jun_17 / Jun_15 / Jun_13 / Jun_10 columns. One analyst
# (row 11) picked up a downgrade alert that week, highlighted in orange.
# Two new research houses (Benchmark, Evercore ISI) were also added as new
# rows at the bottom of the analyst list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 27

# 1) Make room for the new "Jun_27" / "Jun_26" / "Jun_26" group by
#    inserting three columns before the existing first data column (B).
#    Excel shifts the old B/C/D/E (Jun_17/Jun_15/Jun_13/Jun_10) data and
#    their per-cell styles right to E/F/G/H automatically.
$ws.Range("B1:D1").EntireColumn.Insert()

# 2) New header labels for the inserted group.
$ws.Range("B1").Value2 = "Jun_27"
$ws.Range("C1").Value2 = "Jun_26"
$ws.Range("D1").Value2 = "Jun_26"

# 3) Fill the new columns with the same "unchanged" placeholder ("UN")
#    used throughout the rest of the grid for every analyst row.
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $ws.Cells.Item($r, 2).Value2 = "UN"
    $ws.Cells.Item($r, 3).Value2 = "UN"
    $ws.Cells.Item($r, 4).Value2 = "UN"
}

# 4) Row 11 (ValuEngine) got downgraded on 6/21/2018 - flag it in both new
#    "Jun_26" columns with the alert text and an orange highlight, the
#    same way the existing alert cells (rows 3, 6, 20) are highlighted.
$ws.Range("C11").Value2 = "6/21/2018,Downgrades,Hold -> Sell,"
$ws.Range("D11").Value2 = "6/21/2018,Downgrades,Hold -> Sell,"
$ws.Range("C11:D11").Interior.ColorIndex = 45

# 5) Two new research houses added to the bottom of the analyst list.
$ws.Range("A28").Value2 = "Benchmark"
$ws.Range("B28").Value2 = "UN"
$ws.Range("C28").Value2 = "UN"
$ws.Range("D28").Value2 = "UN"

$ws.Range("A29").Value2 = "Evercore ISI"
$ws.Range("B29").Value2 = "UN"
$ws.Range("C29").Value2 = "UN"
$ws.Range("D29").Value2 = "UN"
